# Generate Report for Handback
# Adds a new handed-back file (65d6a1ed-475b-4df5-9188-cdb33c80c9c7.md) as a
# new row on the Overview sheet and on each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newFile      = "65d6a1ed-475b-4df5-9188-cdb33c80c9c7.md"
$newFilePath  = "e2e\65d6a1ed-475b-4df5-9188-cdb33c80c9c7.md"
$newExt       = ".md"
$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newFilePath
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6086e0ff5ff3448452adfcddf64857c4ac609f75/e2e/$newFile", "", "", $newFilePath) | Out-Null
$wsOverview.Range("C3").Value = $newExt
$wsOverview.Range("E3").Value = $statusInSync
$wsOverview.Range("F3").Value = $statusInSync
$wsOverview.Range("G3").Value = "2016-11-29 02:37:17"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6086e0ff5ff3448452adfcddf64857c4ac609f75/e2e/$newFile", "", "", $newFile) | Out-Null
$wsZhCn.Range("B3").Value = $newExt
$wsZhCn.Range("C3").Value = $statusInSync
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "65d6a1ed-475b-4df5-9188-cdb33c80c9c7.83507f8e764a66b020b68463022952d33ad16c8a.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-11-29 02:37:04"
$wsZhCn.Range("I3").Value = $newFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9c36c443af38044d612e1a1a938eeac951ca9ce5/e2e/$newFile", "", "", $newFile) | Out-Null
$wsZhCn.Range("J3").Value = "65d6a1ed-475b-4df5-9188-cdb33c80c9c7.83507f8e764a66b020b68463022952d33ad16c8a.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-11-29 02:37:57"
$wsZhCn.Cells.Item(3, 12).Value = "'"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Cells.Item(3, 14).Value = "'"
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Cells.Item(3, 16).Value = "'"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6086e0ff5ff3448452adfcddf64857c4ac609f75/e2e/$newFile", "", "", $newFile) | Out-Null
$wsDeDe.Range("B3").Value = $newExt
$wsDeDe.Range("C3").Value = $statusInSync
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "65d6a1ed-475b-4df5-9188-cdb33c80c9c7.83507f8e764a66b020b68463022952d33ad16c8a.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-11-29 02:37:17"
$wsDeDe.Range("I3").Value = $newFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/21cef986964300aa1db9d488f963179d004acc36/e2e/$newFile", "", "", $newFile) | Out-Null
$wsDeDe.Range("J3").Value = "65d6a1ed-475b-4df5-9188-cdb33c80c9c7.83507f8e764a66b020b68463022952d33ad16c8a.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-11-29 02:38:15"
$wsDeDe.Cells.Item(3, 12).Value = "'"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Cells.Item(3, 14).Value = "'"
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Cells.Item(3, 16).Value = "'"
